$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-30 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-01 Monday", 2) | Out-Null
$d.Content.Find.Execute("80×68=5440", $true, $false, $false, $false, $false, $true, 1, $false, "53×23=1219", 2) | Out-Null
$d.Content.Find.Execute("77×88=6776", $true, $false, $false, $false, $false, $true, 1, $false, "95×70=6650", 2) | Out-Null
$d.Content.Find.Execute("50×79=3950", $true, $false, $false, $false, $false, $true, 1, $false, "58×64=3712", 2) | Out-Null
$d.Content.Find.Execute("75×58=4350", $true, $false, $false, $false, $false, $true, 1, $false, "86×44=3784", 2) | Out-Null
$d.Content.Find.Execute("87×12=1044", $true, $false, $false, $false, $false, $true, 1, $false, "85×17=1445", 2) | Out-Null
$d.Content.Find.Execute("62×80=4960", $true, $false, $false, $false, $false, $true, 1, $false, "25×87=2175", 2) | Out-Null
$d.Content.Find.Execute("91×70=6370", $true, $false, $false, $false, $false, $true, 1, $false, "33×55=1815", 2) | Out-Null
$d.Content.Find.Execute("11×92=1012", $true, $false, $false, $false, $false, $true, 1, $false, "99×18=1782", 2) | Out-Null
$d.Content.Find.Execute("75×75=5625", $true, $false, $false, $false, $false, $true, 1, $false, "93×99=9207", 2) | Out-Null
$d.Content.Find.Execute("21×61=1281", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=570", 2) | Out-Null
$d.Content.Find.Execute("83×62=5146", $true, $false, $false, $false, $false, $true, 1, $false, "75×30=2250", 2) | Out-Null
$d.Content.Find.Execute("65×42=2730", $true, $false, $false, $false, $false, $true, 1, $false, "46×83=3818", 2) | Out-Null
$d.Content.Find.Execute("84×32=2688", $true, $false, $false, $false, $false, $true, 1, $false, "75×63=4725", 2) | Out-Null
$d.Content.Find.Execute("19×80=1520", $true, $false, $false, $false, $false, $true, 1, $false, "59×60=3540", 2) | Out-Null
$d.Content.Find.Execute("96×60=5760", $true, $false, $false, $false, $false, $true, 1, $false, "70×82=5740", 2) | Out-Null
$d.Content.Find.Execute("37×14=518", $true, $false, $false, $false, $false, $true, 1, $false, "18×57=1026", 2) | Out-Null
$d.Content.Find.Execute("23×77=1771", $true, $false, $false, $false, $false, $true, 1, $false, "61×76=4636", 2) | Out-Null
$d.Content.Find.Execute("56×37=2072", $true, $false, $false, $false, $false, $true, 1, $false, "28×17=476", 2) | Out-Null
$d.Content.Find.Execute("24×20=480", $true, $false, $false, $false, $false, $true, 1, $false, "40×51=2040", 2) | Out-Null
$d.Content.Find.Execute("14×37=518", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=4026", 2) | Out-Null
$d.Content.Find.Execute("98×51=4998", $true, $false, $false, $false, $false, $true, 1, $false, "58×40=2320", 2) | Out-Null
$d.Content.Find.Execute("19×29=551", $true, $false, $false, $false, $false, $true, 1, $false, "79×70=5530", 2) | Out-Null
$d.Content.Find.Execute("89×19=1691", $true, $false, $false, $false, $false, $true, 1, $false, "36×79=2844", 2) | Out-Null
$d.Content.Find.Execute("26×81=2106", $true, $false, $false, $false, $false, $true, 1, $false, "52×49=2548", 2) | Out-Null
$d.Content.Find.Execute("67×24=1608", $true, $false, $false, $false, $false, $true, 1, $false, "13×54=702", 2) | Out-Null
